$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add SENS_PIN label and pin count to row 5 (new encoder sensor pin entry)
$ws.Range("C5").Value = "SENS_PIN"
$ws.Range("D5").Value = 3

# Update running lights pin range from 6~8 to 6~7
$ws.Range("H9").Value = "6~7"

# Update the active selection to D8 as in the saved file
$ws.Range("D8").Select()
